$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.044.48'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").Value = '3.615.49'
$ws.Range("E3").Value = '  +3.39%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.51'
$ws.Range("E5").Value = '  +0.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '195.51'
$ws.Range("E6").Value = '  -0.74%  '
$ws.Range("E7").Value = '  +0.38%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -1.55%  '
$ws.Range("E10").Value = '  -0.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.99'
$ws.Range("E11").Value = '  -0.31%  '
$ws.Range("E12").Value = '  +1.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.55'
$ws.Range("E13").Value = '  -0.13%  '
$ws.Range("D14").Value = '4.188.82'
$ws.Range("E14").Value = '  +3.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.29'
$ws.Range("E15").Value = '  +5.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '592.58'
$ws.Range("E16").Value = '  -1.61%  '
$ws.Range("E17").Value = '  +0.80%  '
$ws.Range("D18").Value = '70.225.69'
$ws.Range("E18").Value = '  +0.43%  '
$ws.Range("D19").Value = '3.614.96'
$ws.Range("E19").Value = '  +3.91%  '
$ws.Range("E20").Value = '  +1.58%  '
$ws.Range("E21").Value = '  +0.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.80'
$ws.Range("E22").Value = '  -2.43%  '
$ws.Range("E23").Value = '  +2.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '102.64'
$ws.Range("E24").Value = '  -1.72%  '
$ws.Range("E25").Value = '  +1.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.05'
$ws.Range("E26").Value = '  -1.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.76'
$ws.Range("E27").Value = '  -2.00%  '
$ws.Range("E28").Value = '  -1.42%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.89'
$ws.Range("E30").Value = '  -0.88%  '
$ws.Range("E31").Value = '  -1.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.34'
$ws.Range("E32").Value = '  -2.78%  '
$ws.Range("E33").Value = '  +1.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.22'
$ws.Range("E34").Value = '  -0.65%  '
$ws.Range("D35").Value = '0.0₃0900'
$ws.Range("E35").Value = '  +11.92%  '
$ws.Range("D36").Value = '3.954.56'
$ws.Range("E36").Value = '  +5.79%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.16'
$ws.Range("E37").Value = '  +5.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '529.22'
$ws.Range("E38").Value = '  +3.98%  '
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.26'
$ws.Range("E40").Value = '  +1.60%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.393'
$ws.Range("E41").Value = '  +0.58%  '
$ws.Range("E42").Value = '  +0.90%  '
$ws.Range("E43").Value = '  -2.04%  '
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("B45").Value = 'ThetaToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.87'
$ws.Range("E45").Value = '  +1.73%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.37'
$ws.Range("E46").Value = '  +1.78%  '
$ws.Range("E47").Value = '  +0.89%  '
$ws.Range("E48").Value = '  -1.36%  '
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("E50").Value = '  +5.36%  '
$ws.Range("E51").Value = '  +3.46%  '

